# Carnet de bord Solomas Ebenisterie
# "implemented a generic function to save ranges + applied it to Appel de fond"
#
# This script reproduces, cell by cell, the edits shown in the author's
# diff for the "Appel de fond" (rows 131-176) and "Eléments de la facture
# à sauvegarder" (rows 182-208) blocks of Feuil1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# --- Row 136/137: relabel the last two "Appel de fond" sub-items ---------
# F136 was "Enregistrement des infos clients" -> becomes the new
# "Enregistrement de la typologie client" entry.
$ws.Range("F136").Value = "Enregistrement de la typologie client"

# F137 was the text "2h" -> becomes a real duration value (0.5 day),
# keeping its existing time-style (s=5) formatting.
$ws.Range("F137").Value = 0.5

# G137 was "08h30" -> becomes "7h45" (still styled s=5).
$ws.Range("G137").Value = "7h45"

# --- Row 143: "Enregistrement des infos clients" block -------------------
# L143 used to be a SUM() formula totalling an (empty) sub-block; it is
# replaced by a plain logged duration, and K143 gets the task label.
$ws.Range("L143").Value = 2
$ws.Range("L143").Style = "Normal"
$ws.Range("K143").Value = "Enregistrement des infos clients"

# --- Row 149: new "Enregistrement des détails de la facture" entry -------
$ws.Range("K149").Value = "Enregistrement des détails de la facture"
$ws.Range("L149").Value = 2

# --- Row 173: "Enregistrement des devis et DMPs" block --------------------
# L173 used to be a SUM() formula totalling an (empty) sub-block; it is
# replaced by a plain logged duration, and K173 gets the task label.
$ws.Range("L173").Value = 2
$ws.Range("K173").Value = "Enregistrement des devis et DMPs"

# --- Row 174: new line describing the named-ranges update on insert ------
$ws.Range("F174").Value = "Insertion d'une nouvelle facture: mise à jour des ranges nommés"
$ws.Range("K174").Value = "Enregistrement de la typologie client"
$ws.Range("L174").Value = 0.5

# --- Row 175: new line describing the named-ranges update on export ------
$ws.Range("F175").Value = 0.5
$ws.Range("F175").NumberFormat = $ws.Range("F137").NumberFormat
$ws.Range("G175").Value = "7h45"
$ws.Range("G175").NumberFormat = $ws.Range("G137").NumberFormat
$ws.Range("K175").Value = "Export d'une nouvelle facture: mise à jour des ranges nommés"
$ws.Range("L175").Value = 0.5

# --- Row 176: new line for the generic "Appel de fond" export function ---
$ws.Range("K176").Value = "Export appel de fond avec ne fonction générique"
$ws.Range("L176").Value = 3

# --- Rows 182-191: "Eléments de la facture à sauvegarder" section --------
# Rows 184 & 185 are repurposed: their old text labels move down to make
# room for two freshly logged work items (K/L pairs), matching the layout
# already used by rows 182/183.
$ws.Range("K184").Value = "Export d'une nouvelle facture: mise à jour des ranges nommés"
$ws.Range("L184").Value = 0.5

$ws.Range("K185").Value = "Export appel de fond avec ne fonction générique"
$ws.Range("L185").Value = 3

# The old label list that used to live in L184:L191 is cleared here; it is
# rewritten further down (rows 197-204) preserving its original order.
$ws.Range("L186:L190").Clear()
$ws.Range("L191").Clear()

# --- Rows 197-204: relocated label list -----------------------------------
$ws.Range("L197").Value = "Eléments de la facture à sauvegader"
$ws.Range("L198").Value = "Montant facture HT"
$ws.Range("L199").Value = "Montants TVA N, R, A"
$ws.Range("L200").Value = "Factures d'acompte"
$ws.Range("L201").Value = "Montant restant dû"
$ws.Range("L202").Value = "Montant du marché non encore appelé"

$ws.Range("L203").Value = "Montant HT par typologie de client"
$ws.Range("L203").Font.Color = 255

$ws.Range("L204").Value = "Typologie Client"

# --- View state: restore where the author had scrolled/selected ----------
$ws.Range("K178").Select()
$excel.ActiveWindow.ScrollRow = 168
$excel.ActiveWindow.ScrollColumn = 7
